$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Enemies")

# Update balancing values (Roketeer row 9, Ironclad row 11)
$ws.Range("C9").Value = 20
$ws.Range("C11").Value = 150
$ws.Range("D11").Value = 0.3

# Update the active selection / cell cursor
$ws.Range("E12").Select()
